# Enhance logging and error handling for XTTS v2 model initialization in TTSWorker
# Update the test-results worksheet: new "Femme"/"Homme" labels replace the old
# generic "femme"/"unknown" values, add an "Own voice (français)" option for
# XTTS v2, and mark the Speedy-Speech test as OK.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-6 and 11: previously generic "femme" -> capitalized "Femme"
$ws.Range("F4").Value = "Femme"
$ws.Range("F5").Value = "Femme"
$ws.Range("F6").Value = "Femme"
$ws.Range("F11").Value = "Femme"

# Row 7 (Speedy-Speech): now succeeds and voice gender is known
$ws.Range("E7").Value = "OK"
$ws.Range("F7").Value = "Femme"

# Row 8 (Neural HMM): voice gender now known
$ws.Range("F8").Value = "Femme"

# Row 9 (XTTS v2): now uses the new "Own voice (français)" option
$ws.Range("F9").Value = "Own voice (français)"

# Row 13 (VCTK_p232, homme)
$ws.Range("F13").Value = "Homme"

# Rows 14-15 (VCTK_p273 / VCTK_p278, femme)
$ws.Range("F14").Value = "Femme"
$ws.Range("F15").Value = "Femme"

# Row 16 (VCTK_p279, homme)
$ws.Range("F16").Value = "Homme"

# Row 17 (VCTK_p304, femme)
$ws.Range("F17").Value = "Femme"

# Widen column F to fit the new, longer values (e.g. "Own voice (français)")
$ws.Columns("F").ColumnWidth = 18.45

# Move the active selection to reflect where the author was last working
$ws.Range("F12").Select() | Out-Null
